# Logged Week 15 and simulated Week 16
# Appends this week's per-play/per-game samples to the running log strings
# on YDS/ST, and updates the season-total numeric cells on OFF/DEF/ST/TURNS/PEN.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Append newly logged samples to the running (space separated) data logs.
# ---------------------------------------------------------------------------

$ydsWs = $wb.Worksheets.Item("YDS")

$ydsAppends = @{
    "B2" = "8 -1 4 -3 5 -2 1 1 9 2 2 6 6 3 -1 3 8 0 0 7 3 -1 4 1 1 -1 3 -3 1 0 -5"
    "C2" = "3 6 2 3 5 2 -1 6 6 3 8 0 3 12 30 2 3 2 6 0 17"
    "B3" = "8 9 7 40 6 7 33 8 4 3 17 5 7"
    "C3" = "5 7 7 -2 3 23 5 8 1 5 5 6 14 19 7 5 8 5 16 6 10 3 6 18 11 13"
}

foreach ($addr in $ydsAppends.Keys) {
    $cell = $ydsWs.Range($addr)
    $existing = $cell.Value()
    $cell.Value = $existing + " " + $ydsAppends[$addr]
}

$stWs = $wb.Worksheets.Item("ST")

$stAppends = @{
    "B4" = "65 60"
    "B5" = "17 22"
    "D3" = "45 58 33 46 54 48 46 49 33"
    "D4" = "6 11 3 2 9 10 0 0 0"
    "D5" = "0 34 0 0 4 8 0 0"
}

foreach ($addr in $stAppends.Keys) {
    $cell = $stWs.Range($addr)
    $existing = $cell.Value()
    $cell.Value = $existing + " " + $stAppends[$addr]
}

# ---------------------------------------------------------------------------
# 2) Update season-total numeric cells that moved as a result of this week's
#    logged/simulated games.
# ---------------------------------------------------------------------------

$offWs = $wb.Worksheets.Item("OFF")
$offValues = @{
    "C2" = 435; "E2" = 22; "F2" = 139; "G2" = 147; "H2" = 11; "I2" = 16; "J2" = 75
    "L2" = 547; "M2" = 350; "Q2" = 1062
    "B3" = 18; "C3" = 329; "E3" = 61; "F3" = 186; "G3" = 58; "H3" = 57; "I3" = 112
    "J3" = 97; "N3" = 22
}
foreach ($addr in $offValues.Keys) {
    $offWs.Range($addr).Value = $offValues[$addr]
}

$defWs = $wb.Worksheets.Item("DEF")
$defValues = @{
    "C2" = 370; "F2" = 105; "G2" = 98; "I2" = 10; "J2" = 61
    "L2" = 609; "M2" = 383; "O2" = 47; "Q2" = 1026
    "C3" = 341; "E3" = 67; "F3" = 220; "G3" = 62; "I3" = 131; "J3" = 116; "N3" = 39
}
foreach ($addr in $defValues.Keys) {
    $defWs.Range($addr).Value = $defValues[$addr]
}

$stValues = @{
    "B2" = 164; "D2" = 133; "J2" = 42; "K2" = 38; "L2" = 18; "M2" = 10
    "B3" = 117
}
foreach ($addr in $stValues.Keys) {
    $stWs.Range($addr).Value = $stValues[$addr]
}

$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsValues = @{
    "C2" = 16; "D2" = 14; "E2" = 10
    "D3" = 13
}
foreach ($addr in $turnsValues.Keys) {
    $turnsWs.Range($addr).Value = $turnsValues[$addr]
}

$penWs = $wb.Worksheets.Item("PEN")
$penValues = @{
    "B2" = 27; "B3" = 27; "D4" = 26
}
foreach ($addr in $penValues.Keys) {
    $penWs.Range($addr).Value = $penValues[$addr]
}
